# Updated cryptos list on Thu Sep  7 06:42:40 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates: values that Excel will not misinterpret as numbers
# (they contain multiple "." separators, or are percent strings padded
# with spaces) can be written directly.
$textUpdates = @(
    @{ Cell = "D2";  Value = "25.805.99" }
    @{ Cell = "E2";  Value = "  -0.15%  " }
    @{ Cell = "D3";  Value = "1.639.67" }
    @{ Cell = "E4";  Value = "  -0.17%  " }
    @{ Cell = "E5";  Value = "  +0.51%  " }
    @{ Cell = "E7";  Value = "  -0.15%  " }
    @{ Cell = "E8";  Value = "  -0.01%  " }
    @{ Cell = "E9";  Value = "  -0.91%  " }
    @{ Cell = "E10"; Value = "  -0.88%  " }
    @{ Cell = "E11"; Value = "  +1.52%  " }
    @{ Cell = "E12"; Value = "  +0.32%  " }
    @{ Cell = "D13"; Value = "1.865.87" }
    @{ Cell = "E13"; Value = "  +0.28%  " }
    @{ Cell = "D14"; Value = "1.638.63" }
    @{ Cell = "E14"; Value = "  +0.24%  " }
    @{ Cell = "E15"; Value = "  +0.48%  " }
    @{ Cell = "E16"; Value = "  +0.16%  " }
    @{ Cell = "E17"; Value = "  +0.01%  " }
    @{ Cell = "D18"; Value = "25.845.96" }
    @{ Cell = "E19"; Value = "  -0.13%  " }
    @{ Cell = "E20"; Value = "  +2.33%  " }
    @{ Cell = "E21"; Value = "  -0.43%  " }
    @{ Cell = "E22"; Value = "  +0.74%  " }
    @{ Cell = "E23"; Value = "  +2.22%  " }
    @{ Cell = "E24"; Value = "  +4.69%  " }
    @{ Cell = "E25"; Value = "  -0.12%  " }
    @{ Cell = "E26"; Value = "  +2.19%  " }
    @{ Cell = "E28"; Value = "  +1.82%  " }
    @{ Cell = "E29"; Value = "  -0.03%  " }
    @{ Cell = "E30"; Value = "  +0.29%  " }
    @{ Cell = "E31"; Value = "  -0.13%  " }
    @{ Cell = "E32"; Value = "  +0.86%  " }
    @{ Cell = "E33"; Value = "  -0.25%  " }
    @{ Cell = "E34"; Value = "  +0.65%  " }
    @{ Cell = "E35"; Value = "  -0.15%  " }
    @{ Cell = "D37"; Value = "1.133.16" }
    @{ Cell = "E37"; Value = "  +0.94%  " }
    @{ Cell = "E38"; Value = "  -1.88%  " }
    @{ Cell = "E39"; Value = "  -0.33%  " }
    @{ Cell = "E40"; Value = "  +0.42%  " }
    @{ Cell = "E41"; Value = "  -0.07%  " }
    @{ Cell = "E42"; Value = "  +1.48%  " }
    @{ Cell = "E43"; Value = "  +1.25%  " }
    @{ Cell = "E44"; Value = "  +0.78%  " }
    @{ Cell = "D45"; Value = "1.775.06" }
    @{ Cell = "E45"; Value = "  +0.24%  " }
    @{ Cell = "E46"; Value = "  +3.16%  " }
    @{ Cell = "E47"; Value = "  -0.03%  " }
    @{ Cell = "E48"; Value = "  -1.29%  " }
    @{ Cell = "E49"; Value = "  +5.12%  " }
    @{ Cell = "E51"; Value = "  +1.81%  " }
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Price updates whose new text looks like a plain decimal number
# (e.g. "215.92"). A bare .Value assignment would let Excel auto-convert
# these into numeric cells, which would both change the cell type and can
# silently drop meaningful trailing zeros (e.g. "192.60" -> 192.6).
# Force the cell to Text format first, assign the literal string, then
# restore the original ("Normal") cell style so no formatting residue is
# left behind.
$numericLookingUpdates = @(
    @{ Cell = "D5";  Value = "215.92" }
    @{ Cell = "D10"; Value = "19.68" }
    @{ Cell = "D11"; Value = "0.0793" }
    @{ Cell = "D15"; Value = "0.562" }
    @{ Cell = "D17"; Value = "63.05" }
    @{ Cell = "D20"; Value = "4.49" }
    @{ Cell = "D21"; Value = "192.60" }
    @{ Cell = "D23"; Value = "6.33" }
    @{ Cell = "D26"; Value = "142.19" }
    @{ Cell = "D31"; Value = "0.0493" }
    @{ Cell = "D32"; Value = "3.33" }
    @{ Cell = "D34"; Value = "1.58" }
    @{ Cell = "D39"; Value = "0.546" }
    @{ Cell = "D43"; Value = "100.79" }
    @{ Cell = "D44"; Value = "0.806" }
    @{ Cell = "D47"; Value = "55.37" }
    @{ Cell = "D48"; Value = "0.417" }
    @{ Cell = "D49"; Value = "1.43" }
)

foreach ($u in $numericLookingUpdates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
